# Update TPM (transcripts per million) based NATMI LR-pair metrics on
# Sheet1 (Adm-Ramp2.xlsx) to reflect the newly recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 15.89577633333333
$ws.Range("H2").Value = 47.687329
$ws.Range("I2").Value = 0.286059172443548
$ws.Range("J2").Value = 0.2860591724435479
$ws.Range("M2").Value = 54.11233133333334
$ws.Range("N2").Value = 162.336994
$ws.Range("O2").Value = 0.4576967561138398
$ws.Range("P2").Value = 0.4576967561138398
$ws.Range("Q2").Value = 860.1575157498918
$ws.Range("R2").Value = 7741.417641749026
$ws.Range("S2").Value = 0.1309283552840214
$ws.Range("T2").Value = 0.1309283552840214
$ws.Range("G3").Value = 15.89577633333333
$ws.Range("H3").Value = 47.687329
$ws.Range("I3").Value = 0.286059172443548
$ws.Range("J3").Value = 0.2860591724435479
$ws.Range("O3").Value = 0.5177793530461455
$ws.Range("P3").Value = 0.5177793530461455
$ws.Range("Q3").Value = 973.0717905983682
$ws.Range("R3").Value = 8757.646115385314
$ws.Range("S3").Value = 0.1481155332407361
$ws.Range("T3").Value = 0.148115533240736
$ws.Range("G4").Value = 15.89577633333333
$ws.Range("H4").Value = 47.687329
$ws.Range("I4").Value = 0.286059172443548
$ws.Range("J4").Value = 0.2860591724435479
$ws.Range("O4").Value = 0.02452389084001462
$ws.Range("P4").Value = 0.02452389084001462
$ws.Range("Q4").Value = 46.08817681072189
$ws.Range("R4").Value = 414.793591296497
$ws.Range("S4").Value = 0.007015283918790489
$ws.Range("T4").Value = 0.007015283918790487
$ws.Range("I5").Value = 0.6735478078679881
$ws.Range("J5").Value = 0.673547807867988
$ws.Range("M5").Value = 54.11233133333334
$ws.Range("N5").Value = 162.336994
$ws.Range("O5").Value = 0.4576967561138398
$ws.Range("P5").Value = 0.4576967561138398
$ws.Range("Q5").Value = 2025.305478602846
$ws.Range("R5").Value = 18227.74930742561
$ws.Range("S5").Value = 0.3082806467487659
$ws.Range("T5").Value = 0.3082806467487659
$ws.Range("I6").Value = 0.6735478078679881
$ws.Range("J6").Value = 0.673547807867988
$ws.Range("O6").Value = 0.5177793530461455
$ws.Range("P6").Value = 0.5177793530461455
$ws.Range("S6").Value = 0.3487491482035364
$ws.Range("T6").Value = 0.3487491482035363
$ws.Range("I7").Value = 0.6735478078679881
$ws.Range("J7").Value = 0.673547807867988
$ws.Range("O7").Value = 0.02452389084001462
$ws.Range("P7").Value = 0.02452389084001462
$ws.Range("R7").Value = 976.662666499814
$ws.Range("S7").Value = 0.01651801291568568
$ws.Range("T7").Value = 0.01651801291568568
$ws.Range("I8").Value = 0.04039301968846393
$ws.Range("J8").Value = 0.04039301968846393
$ws.Range("M8").Value = 54.11233133333334
$ws.Range("N8").Value = 162.336994
$ws.Range("O8").Value = 0.4576967561138398
$ws.Range("P8").Value = 0.4576967561138398
$ws.Range("Q8").Value = 121.4586449792034
$ws.Range("R8").Value = 1093.12780481283
$ws.Range("S8").Value = 0.0184877540810524
$ws.Range("T8").Value = 0.0184877540810524
$ws.Range("I9").Value = 0.04039301968846393
$ws.Range("J9").Value = 0.04039301968846393
$ws.Range("O9").Value = 0.5177793530461455
$ws.Range("P9").Value = 0.5177793530461455
$ws.Range("S9").Value = 0.02091467160187307
$ws.Range("T9").Value = 0.02091467160187307
$ws.Range("I10").Value = 0.04039301968846393
$ws.Range("J10").Value = 0.04039301968846393
$ws.Range("O10").Value = 0.02452389084001462
$ws.Range("P10").Value = 0.02452389084001462
$ws.Range("S10").Value = 0.0009905940055384508
$ws.Range("T10").Value = 0.0009905940055384506
